$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 338.5
$ws.Range("B2").Value = 291.5
$ws.Range("C2").Value = 16.081743089950155
$ws.Range("D2").Value = 34.897565432435158
$ws.Range("E2").Value = 15.733901658773538
$ws.Range("F2").Value = 30.940781288051532
